# Update gh-pages to output generated at 456a3b4
#
# The upstream scraper re-ran and produced a fresh export: the "开始时间"
# (start time) column now renders dates with dots instead of dashes
# (e.g. "2024-02-15" -> "2024.02.15"), and the "想去人数" (want-to-go
# count) column picked up newer crawl numbers for a few rows. Apply both
# across the three sheets that carry these rows ("展览", "演出",
# "全部类型" - "本地生活" has no data rows and is untouched).
#
# Note: the date-like strings must stay plain text (they are stored as
# inlineStr/shared-string cells, not real dates). Setting .Value directly
# on a string like "2024.02.15" makes Excel auto-convert it to a date
# serial, so we briefly force Text number-formatting for the write and
# then restore the cell's style back to Normal so no stray formatting is
# left behind.

function Set-TextValue {
    param($Range, [string]$Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
Set-TextValue $ws.Range("B2") "2024.02.15"
$ws.Range("F2").Value = 10244
Set-TextValue $ws.Range("B3") "2024.03.09"
$ws.Range("F3").Value = 229
Set-TextValue $ws.Range("B4") "2024.03.16"
$ws.Range("F4").Value = 55
Set-TextValue $ws.Range("B5") "2024.03.16"
$ws.Range("F5").Value = 637
Set-TextValue $ws.Range("B6") "2024.03.30"
$ws.Range("F6").Value = 484

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
Set-TextValue $ws.Range("B2") "2024.03.30"

# --- Sheet "本地生活" : no changes ---

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
Set-TextValue $ws.Range("B2") "2024.02.15"
$ws.Range("F2").Value = 10244
Set-TextValue $ws.Range("B3") "2024.03.09"
$ws.Range("F3").Value = 229
Set-TextValue $ws.Range("B4") "2024.03.16"
$ws.Range("F4").Value = 55
Set-TextValue $ws.Range("B5") "2024.03.16"
$ws.Range("F5").Value = 637
Set-TextValue $ws.Range("B6") "2024.03.30"
Set-TextValue $ws.Range("B7") "2024.03.30"
$ws.Range("F7").Value = 484
